$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.622.33"
$ws.Range("E2").Value = "  -1.68%  "
$ws.Range("D3").Value = "1.966.90"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("D4").Value = "'1.011"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'323.55"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.24%  "
$ws.Range("D6").Value = "'1.011"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("D7").Value = "'0.4831"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -3.12%  "
$ws.Range("D8").Value = "'0.4079"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.29%  "
$ws.Range("D9").Value = "'54.17"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.30%  "
$ws.Range("D10").Value = "'0.08521"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -5.76%  "
$ws.Range("E11").Value = "  -3.22%  "
$ws.Range("D12").Value = "'22.53"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.37%  "
$ws.Range("D13").Value = "2.002.78"
$ws.Range("E13").Value = "  +4.08%  "
$ws.Range("D14").Value = "'7.635"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -3.16%  "
$ws.Range("D15").Value = "'6.212"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.65%  "
$ws.Range("E16").Value = "  +0.45%  "
$ws.Range("D17").Value = "'91.19"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("D18").Value = "'0.00001076"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -2.15%  "
$ws.Range("D19").Value = "'0.06641"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.72%  "
$ws.Range("D20").Value = "'18.64"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.83%  "
$ws.Range("E21").Value = "  +0.47%  "
$ws.Range("D22").Value = "'5.896"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.91%  "
$ws.Range("D23").Value = "28.663.08"
$ws.Range("E23").Value = "  -1.55%  "
$ws.Range("E24").Value = "  -3.25%  "
$ws.Range("D25").Value = "'2.303"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.53%  "
$ws.Range("D26").Value = "2.244.26"
$ws.Range("E26").Value = "  +3.08%  "
$ws.Range("D27").Value = "'156.74"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.52%  "
$ws.Range("D28").Value = "'20.41"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.09%  "
$ws.Range("E29").Value = "  -3.63%  "
$ws.Range("D30").Value = "'2.194"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.01%  "
$ws.Range("D31").Value = "'125.02"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.85%  "
$ws.Range("D32").Value = "'0.9973"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -4.48%  "
$ws.Range("D33").Value = "'0.09715"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.31%  "
$ws.Range("D34").Value = "'1.472"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -4.30%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").Value = "'5.662"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.36%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "'3.701"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.34%  "
$ws.Range("D37").Value = "'9.194"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.85%  "
$ws.Range("D38").Value = "'0.02347"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.18%  "
$ws.Range("D39").Value = "'0.06271"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.69%  "
$ws.Range("D40").Value = "'1.259"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.52%  "
$ws.Range("D41").Value = "'0.6267"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.94%  "
$ws.Range("D42").Value = "'11.29"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.69%  "
$ws.Range("E43").Value = "  +0.44%  "
$ws.Range("D44").Value = "'0.1925"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.37%  "
$ws.Range("E45").Value = "  +6.03%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.5987"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.64%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'13.10"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.61%  "
$ws.Range("D48").Value = "'2.080"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -4.34%  "
$ws.Range("D49").Value = "'3.419"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.45%  "
$ws.Range("D50").Value = "'0.06855"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.56%  "
$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").Value = "'111.89"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.49%  "
